# Staging.Organization.xlsx - "staging templates updated with database changes"
#
# The header row (row 2) columns E:G are reordered from
#   E=ShortName, F=LongName, G=ParentOrganization_ID
# to
#   E=LongName, F=ParentOrganization_ID, G=ShortName
#
# (LongName now comes right after Code, followed by ParentOrganization_ID,
#  then ShortName last - matching the updated database column ordering.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "LongName"
$ws.Range("F2").Value = "ParentOrganization_ID"
$ws.Range("G2").Value = "ShortName"

# Restore default A1 selection (best effort - matches the author dropping the
# explicit <selection> from the saved view state).
$ws.Range("A1").Select()
